# "wrong report picture uploaded - fixed"
#
# The "Метод ветвей и границ с ЛП" (row 23) block had been pasted with the
# wrong picture's numbers. The correct values are identical to the
# "ДП по стоимостям" row (row 21) directly above it in the same table.
# Copy that row's values over the bogus ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correct = $ws.Range("C21:I21")
$wrong   = $ws.Range("C23:I23")
$wrong.Value = $correct.Value()

# Leave the cursor where the author last left it after verifying the fix.
$ws.Range("R31").Select()
